$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Fitness column (C) for rows 2 through 252 from 7293 to 7573
$ws.Range("C2:C252").Value = 7573
